# Update imputed values (KNN result) in columns D and E for several rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.314
$ws.Range("E4").Value = 16.545
$ws.Range("D11").Value = -7.084000000000001
$ws.Range("D12").Value = -7.208
$ws.Range("E14").Value = 16.905
$ws.Range("D15").Value = -8.297000000000001
$ws.Range("E26").Value = 17.054
$ws.Range("D27").Value = -8.514999999999999
$ws.Range("D28").Value = -8.108000000000001
$ws.Range("D31").Value = -7.770000000000001
$ws.Range("E31").Value = 16.916
$ws.Range("D32").Value = -7.945000000000002
$ws.Range("E35").Value = 16.472
$ws.Range("D36").Value = -7.723999999999999
$ws.Range("E37").Value = 16.643
$ws.Range("D38").Value = -7.828
$ws.Range("E39").Value = 16.59
$ws.Range("E40").Value = 16.595
$ws.Range("E45").Value = 16.96
$ws.Range("D46").Value = -8.176000000000002
$ws.Range("E52").Value = 16.768
$ws.Range("D54").Value = -8.229999999999999
$ws.Range("D55").Value = -8.129000000000001
$ws.Range("D56").Value = -8.17
$ws.Range("E57").Value = 16.643
$ws.Range("D67").Value = -7.509
$ws.Range("D69").Value = -7.292
$ws.Range("D72").Value = -7.595000000000001
$ws.Range("D73").Value = -7.969000000000001
$ws.Range("E81").Value = 16.957
$ws.Range("D83").Value = -7.852000000000001
$ws.Range("E83").Value = 16.661
$ws.Range("D86").Value = -7.885999999999998
$ws.Range("D91").Value = -7.141000000000001
$ws.Range("D93").Value = -7.855
$ws.Range("D99").Value = -8.030000000000001
$ws.Range("E100").Value = 16.684
$ws.Range("E102").Value = 16.555
